$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = "Comments"
$ws.Range("J2").Font.Bold = $true

$ws.Range("J3").Value = "UP test comment 1"
$ws.Range("J4").Value = "UP test comment 2"
$ws.Range("J5").Value = "UP test comment 3"
$ws.Range("J6").Value = "UP test comment 4"

$ws.Range("J3:J6").Select() | Out-Null

